# Edit slide 3 ("팀원 역할 분담") of the presentation:
#   - 박주영's role-description textbox drops the "서버 클라이언트" portion so the
#     parenthetical reads "도서 관리 담당" instead of "도서 관리 / 서버 클라이언트 담당".
#   - The small layout nudges that PowerPoint applies to every shape on the
#     slide (a fraction of a point) are reproduced, together with the
#     resize of the edited textbox to its new (narrower) autofit width.
#   - The slide transition speed changes from Medium to Slow.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---------------------------------------------------------------------
# 1) Fix the text in "TextBox 20" (박주영 ...) - remove " /  서버  클라이언트  "
# ---------------------------------------------------------------------
$roleBox = $s.Shapes.Item(10)
$tr = $roleBox.TextFrame.TextRange
# Before: "박주영  ( 도서  관리  /  서버  클라이언트  담당 )"
# After : "박주영  ( 도서  관리 담당 )"
$cut = $tr.Characters(15, 15)
$cut.Delete()

# Paragraph no longer carries a bullet definition.
$tr.ParagraphFormat.Bullet.Type = 0

# Body-text insets / anchoring normalized by PowerPoint on save.
$roleBox.TextFrame.MarginLeft = 7.2
$roleBox.TextFrame.MarginTop = 3.6
$roleBox.TextFrame.MarginRight = 7.2
$roleBox.TextFrame.MarginBottom = 3.6
$roleBox.TextFrame.VerticalAnchor = 1
$roleBox.TextFrame.Orientation = 1

# ---------------------------------------------------------------------
# 2) Re-apply the (sub-point) geometry PowerPoint recomputed for every
#    shape on the slide, including the shrink of the edited textbox.
# ---------------------------------------------------------------------
function Set-ShapeGeometry($shape, $left, $top, $width, $height) {
    $shape.Left = $left
    $shape.Top = $top
    $shape.Width = $width
    $shape.Height = $height
}

Set-ShapeGeometry ($s.Shapes.Item(2))  131.99996948242188 25.999961853027344  242.74996948242188 60.5999641418457    # TextBox 11 (title)
Set-ShapeGeometry ($s.Shapes.Item(3))  100.64996337890625 161.24996948242188 91.99996185302734  82.19996643066406   # 직사각형 13
Set-ShapeGeometry ($s.Shapes.Item(4))  211.64996337890625 161.24996948242188 649.0               82.19996643066406  # 직사각형 14
Set-ShapeGeometry ($s.Shapes.Item(5))  130.0999755859375  175.94996643066406 27.399961471557617  50.89996337890625  # TextBox 15 ("1")
Set-ShapeGeometry ($s.Shapes.Item(6))  231.89996337890625 179.5999755859375  361.89996337890625  46.04996109008789  # TextBox 16 (강현영 ...)
Set-ShapeGeometry ($s.Shapes.Item(7))  100.64996337890625 265.14996337890625 91.99996185302734  82.19996643066406   # 직사각형 17
Set-ShapeGeometry ($s.Shapes.Item(8))  211.64996337890625 265.14996337890625 649.0               82.19996643066406  # 직사각형 18
Set-ShapeGeometry ($s.Shapes.Item(9))  127.79996490478516 280.6999816894531  33.749961853027344  50.89996337890625  # TextBox 19 ("2")
Set-ShapeGeometry ($s.Shapes.Item(10)) 231.89996337890625 281.8499755859375  298.39996337890625  45.999961853027344 # TextBox 20 (박주영 ...) - narrower now
Set-ShapeGeometry ($s.Shapes.Item(11)) 100.64996337890625 369.04998779296875 91.99996185302734  82.19996643066406   # 직사각형 21
Set-ShapeGeometry ($s.Shapes.Item(12)) 211.64996337890625 369.04998779296875 649.0               82.19996643066406  # 직사각형 22
Set-ShapeGeometry ($s.Shapes.Item(13)) 127.79996490478516 384.5999755859375  33.749961853027344  50.89996337890625  # TextBox 23 ("3")
Set-ShapeGeometry ($s.Shapes.Item(14)) 231.89996337890625 385.7499694824219  529.75              46.04996109008789  # TextBox 24 (김종혁 ...)

# ---------------------------------------------------------------------
# 3) Slide transition speed: Medium -> Slow
# ---------------------------------------------------------------------
$s.SlideShowTransition.Speed = 1
